# SENG201 project report — testing-section update
#
# 1. Merge the "Model View Controller ... archetectual pattern ..."
#    run-run-run-run paragraph into a single run (no textual change).
# 2. Bump the coverage percentages in the JUnit Testing paragraph
#    (85->89, 94->96, 98->100, 61->64, 59->60) while reproducing the
#    fine-grained run split the original capture shows.
# 3. Merge the "This would have allowed" run-run-run-run-run paragraph
#    into a single run (no textual change).
# 4. Suppress automatic hyphenation on the Normal / No Spacing styles.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force Word to materialise an explicit (empty) <w:rPr/> on the
# run(s) covering a Range, and to split runs at the Range boundaries,
# by toggling a character property on and back off again.
# ---------------------------------------------------------------------
function Stamp-Range($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

# -----------------------------------------------------------------
# 1) "Application Structure" paragraph: collapse the 4 runs that make
#    up the MVC sentence into one run. The wording itself is unchanged.
# -----------------------------------------------------------------
$mvcText = "The game is structured around the idea of a Model View Controller archetectual pattern as learnt about in lectures. This allows for easy maintenance and testing, which are both essential when creating a project like this where there are constant changes happening."
$d.Content.Find.Execute($mvcText, $true, $false, $false, $false, $false, $true, 1, $false, $mvcText, 2) | Out-Null

# find that paragraph again (index is stable) and stamp it so the
# merged run carries an explicit empty <w:rPr/>
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq $mvcText) {
        $r = $para.Range
        Stamp-Range $d.Range($r.Start, $r.End - 1)
        break
    }
}

# -----------------------------------------------------------------
# 2) "Junit Testing" paragraph: update the coverage percentages.
# -----------------------------------------------------------------
$d.Content.Find.Execute("85% class", $true, $false, $false, $false, $false, $true, 1, $false, "89% class", 2) | Out-Null
$d.Content.Find.Execute("94% method", $true, $false, $false, $false, $false, $true, 1, $false, "96% method", 2) | Out-Null
$d.Content.Find.Execute("98% method", $true, $false, $false, $false, $false, $true, 1, $false, "100% method", 2) | Out-Null
$d.Content.Find.Execute("61% class", $true, $false, $false, $false, $false, $true, 1, $false, "64% class", 2) | Out-Null
$d.Content.Find.Execute("59% method", $true, $false, $false, $false, $false, $true, 1, $false, "60% method", 2) | Out-Null

$coverageText = "In our experience with running our tests we were able to achieve 89% class coverage and 96% method coverage in the models package. This is due to the fact that we aren" + [char]8217 + "t testing all of the classes withing the towertypes package as they are all just classes constructing the starting towers and so they don" + [char]8217 + "t require testing. We also achieved 100% coverage in the service package with 100% method coverage. Overall however, the entire project only had a 64% class coverage with a 60% method coverage. This can easily be explained by the fact that there isn" + [char]8217 + "t a straightforward way to test the controller classes for the GUI and so they remain untested."

$coverageStart = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq $coverageText) {
        $coverageStart = $para.Range.Start
        break
    }
}

# Boundaries (relative to paragraph start) of the 11 runs the capture
# shows: big chunk / "9" / chunk / "6" / chunk / "100" / chunk / "4" /
# chunk / "60" / chunk.
$bounds = @(0, 66, 67, 89, 90, 384, 387, 454, 455, 479, 481, 651)
for ($i = 0; $i -lt ($bounds.Length - 1); $i++) {
    $s = $coverageStart + $bounds[$i]
    $e = $coverageStart + $bounds[$i + 1]
    Stamp-Range ($d.Range($s, $e))
}

# -----------------------------------------------------------------
# 3) "Retrospective" paragraph: collapse the 5 runs that make up the
#    "This would have allowed" sentence into one run. Wording unchanged.
# -----------------------------------------------------------------
$retroText = "Overall, we felt like the project as a whole went quite smoothly. We started as possible and worked throughout the holidays in order to give us as much opportunity to get the game in a state that we were happy with. We were able to implement a lot of our learnings from the course and felt like we learnt a lot from doing that. In future however, it would be essential to allocate more time at the beginning towards planning in order to have a better idea of what we would be building. This would have allowed for a much smoother experience when creating an application as we would have a set idea of what we would have to do."
$d.Content.Find.Execute($retroText, $true, $false, $false, $false, $false, $true, 1, $false, $retroText, 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq $retroText) {
        $r = $para.Range
        Stamp-Range $d.Range($r.Start, $r.End - 1)
        break
    }
}

# -----------------------------------------------------------------
# 4) Suppress automatic hyphenation for the Normal and No Spacing
#    paragraph styles.
# -----------------------------------------------------------------
$d.Styles.Item("Normal").ParagraphFormat.Hyphenation = $false
$d.Styles.Item("No Spacing").ParagraphFormat.Hyphenation = $false
